# Remove the redundant "TextBox 59" shape (the duplicate aspect-label
# textbox) from slide 7 of the presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$s.Shapes.Item("TextBox 59").Delete()
